$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '48.295.21'
$ws.Range('E2').Value = '  +2.33%  '

$ws.Range('D3').Value = '2.510.62'
$ws.Range('E3').Value = '  +1.06%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = "'321.39"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.17%  '

$ws.Range('D6').Value = "'108.73"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.64%  '

$ws.Range('E7').Value = '  +1.23%  '

$ws.Range('D9').Value = "'0.544"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.93%  '

$ws.Range('D10').Value = "'39.95"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.23%  '

$ws.Range('D11').Value = "'20.08"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +9.17%  '

$ws.Range('E12').Value = '  +1.04%  '

$ws.Range('D13').Value = "'0.124"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.24%  '

$ws.Range('D14').Value = "'7.20"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.52%  '

$ws.Range('D15').Value = '2.908.18'
$ws.Range('E15').Value = '  +1.20%  '

$ws.Range('D16').Value = '2.518.72'
$ws.Range('E16').Value = '  +1.21%  '

$ws.Range('D17').Value = "'0.846"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.35%  '

$ws.Range('D18').Value = '48.130.24'
$ws.Range('E18').Value = '  +2.14%  '

$ws.Range('D19').Value = "'13.12"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.59%  '

$ws.Range('D20').Value = "'6.75"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.00%  '

$ws.Range('E21').Value = '  +0.50%  '

$ws.Range('E22').Value = '  +0.10%  '

$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = "'72.27"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.58%  '

$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = "'277.68"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +12.96%  '

$ws.Range('D25').Value = "'2.56"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.70%  '

$ws.Range('E26').Value = '  +0.00%  '

$ws.Range('D27').Value = "'25.84"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.77%  '

$ws.Range('D28').Value = "'2.40"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.81%  '

$ws.Range('D29').Value = "'9.84"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.01%  '

$ws.Range('D30').Value = "'35.48"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.70%  '

$ws.Range('E31').Value = '  -0.07%  '

$ws.Range('D32').Value = "'49.18"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.18%  '

$ws.Range('D33').Value = "'19.54"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.62%  '

$ws.Range('E34').Value = '  +0.76%  '

$ws.Range('D35').Value = "'1.00"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.08%  '

$ws.Range('D36').Value = "'0.0785"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.51%  '

$ws.Range('D37').Value = "'1.96"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.66%  '

$ws.Range('E38').Value = '  -2.46%  '

$ws.Range('D39').Value = "'2.95"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.91%  '

$ws.Range('E40').Value = '  +0.15%  '

$ws.Range('D41').Value = "'122.22"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.66%  '

$ws.Range('D42').Value = "'2.21"
$ws.Range('D42').Style = 'Normal'

$ws.Range('D43').Value = "'21.69"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.43%  '

$ws.Range('E44').Value = '  +3.45%  '

$ws.Range('D45').Value = '2.000.94'
$ws.Range('E45').Value = '  +0.32%  '

$ws.Range('E46').Value = '  +5.21%  '

$ws.Range('E47').Value = '  +3.25%  '

$ws.Range('E48').Value = '  -0.75%  '

$ws.Range('D49').Value = "'9.03"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.64%  '

$ws.Range('E50').Value = '  +2.67%  '

$ws.Range('D51').Value = "'80.31"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.70%  '
